# Corrections made during class review.
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Footer "date" placeholder text: 30/04/2012 -> 7/7/2012
#    Present once in the Slide Master, once in each of the 11
#    Custom Layouts, and once in the Notes Master.
# ---------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "30/04/2012") {
                $shp.TextFrame.TextRange.Text = "7/7/2012"
            }
        }
    }
}

Update-DatePlaceholder($p.SlideMaster)

for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    Update-DatePlaceholder($p.SlideMaster.CustomLayouts.Item($L))
}

Update-DatePlaceholder($p.NotesMaster)

# ---------------------------------------------------------------
# 2) Slide 12 ("Classpath"): merge the run-split of the second
#    paragraph back together (same visible text, cleaner runs).
# ---------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp12 = $s12.Shapes.Item(2)
$tr12 = $shp12.TextFrame.TextRange
$para12 = $tr12.Paragraphs(2, 1)

$restStart = 12
$restLen = $para12.Length - ($restStart - 1)
$chRest = $para12.Characters($restStart, $restLen)
$chRest.Text = " precisa ser especificado tanto ao compilar quanto ao executar suas classes"

$chWord = $para12.Characters(3, 9)
$chWord.Text = "classpath"

# ---------------------------------------------------------------
# 3) Slide 5 ("Como criar um arquivo JAR"): jar -C path changes
#    from C:\src to C:\bin.
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(4, 1)
$ch5 = $para5.Characters(41, 3)
$ch5.Text = "bin"

# ---------------------------------------------------------------
# 4) Slide 8 ("O arquivo MANIFEST.MF"): same C:\src -> C:\bin fix.
# ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(3, 1)
$ch8 = $para8.Characters(51, 3)
$ch8.Text = "bin"
